# Modify SLG building config: insert "Icon" and "ShowName" columns
# (G and H) before the existing "Desc" column, which moves to I.
# Icon values are the bare prefab/object names; ShowName duplicates
# the existing Desc (localized display name) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Cells.Item(1, 7).Value = "Icon"
$ws.Cells.Item(1, 8).Value = "ShowName"
$ws.Cells.Item(1, 9).Value = "Desc"

$rows = @(
    @{ Icon = "Altar_1_1";         Name = "一级祭坛" },
    @{ Icon = "Arena_1_1";         Name = "一级竞技场" },
    @{ Icon = "Camp_1_1";          Name = "一级兵营" },
    @{ Icon = "GoldMine_1_1";      Name = "一级金矿" },
    @{ Icon = "Item_hourse_1_1";   Name = "一级道具屋" },
    @{ Icon = "League_1_1";        Name = "一级公会" },
    @{ Icon = "MagicHourse_1_1";   Name = "一级魔法屋" },
    @{ Icon = "Tower_1_1";         Name = "一级箭塔" },
    @{ Icon = "Town_1_1";          Name = "一级大厅" }
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]
    $ws.Cells.Item($r, 7).Value = $data.Icon
    $ws.Cells.Item($r, 8).Value = $data.Name
    $ws.Cells.Item($r, 9).Value = $data.Name
    # Data rows use the sheet's "Text" number format (same style as the
    # rest of the row); new H/I columns don't inherit it automatically.
    $ws.Range("G" + $r + ":I" + $r).NumberFormat = "@"
}

# Column width: merge old G (11) / H (14) widths into a single
# G:I width-11 block, matching the new layout. (10.2857... compensates
# for Excel's internal pixel-padding round-trip so the stored XML width
# comes out to exactly 11, as in the original G column.)
$ws.Range("G1:I1").EntireColumn.ColumnWidth = 10.285714285714286

# Selection marker left by the editor after the last entry (H10)
$ws.Range("H10").Select()
